$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '28.131.33'
$ws.Cells.Item(2, 5).Value = '  +0.44%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.875.34'
$ws.Cells.Item(3, 5).Value = '  -1.25%  '

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.001'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.27%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '313.43'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.30%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  -0.19%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5078'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  +1.08%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3845'
$ws.Cells.Item(8, 4).ClearFormats()
$ws.Cells.Item(8, 5).Value = '  -2.02%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.09017'
$ws.Cells.Item(9, 4).ClearFormats()
$ws.Cells.Item(9, 5).Value = '  -5.48%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.123'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  -0.74%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '41.60'
$ws.Cells.Item(11, 4).ClearFormats()
$ws.Cells.Item(11, 5).Value = '  -0.77%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '6.356'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  -0.15%  '

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '20.80'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +0.20%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.871.05'
$ws.Cells.Item(14, 5).Value = '  -1.90%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.43%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  -0.26%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  -0.75%  '

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '91.31'
$ws.Cells.Item(18, 4).ClearFormats()
$ws.Cells.Item(18, 5).Value = '  -0.82%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06576'
$ws.Cells.Item(19, 4).ClearFormats()
$ws.Cells.Item(19, 5).Value = '  -0.35%  '

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.31'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +2.88%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '1.0000'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  -0.14%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.139'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  -0.85%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '28.153.16'
$ws.Cells.Item(23, 5).Value = '  +0.18%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.45'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +1.56%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.265'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  -1.57%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(26, 4).Value = '2.094.44'
$ws.Cells.Item(26, 5).Value = '  -1.22%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.540'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  -4.17%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '20.86'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +0.30%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Monero'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '157.18'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +0.09%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'BitcoinCash'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '127.05'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  +0.14%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Stellar'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.1064'
$ws.Cells.Item(31, 4).ClearFormats()
$ws.Cells.Item(31, 5).Value = '  +0.05%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'ImmutableX'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.059'
$ws.Cells.Item(32, 4).ClearFormats()
$ws.Cells.Item(32, 5).Value = '  -1.81%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.617'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +0.17%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'HuobiToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.604'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  -0.51%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'FraxShare'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '9.499'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  -0.42%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'Hedera'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.06597'
$ws.Cells.Item(36, 4).ClearFormats()
$ws.Cells.Item(36, 5).Value = '  +0.13%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.02401'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  -1.09%  '

# Row 38
$ws.Cells.Item(38, 2).Value = 'Algorand'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2199'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +1.12%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.297'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +2.68%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'ARBITRUM'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.212'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  -1.39%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'TheSandbox'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6433'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +1.80%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Aptos'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '11.53'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  +1.89%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '4.931'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  -1.07%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Decentraland'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.6038'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +1.05%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'EnergySwap'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '13.23'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  -0.65%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.276'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -0.15%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'PancakeSwap'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.666'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -1.47%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'EOS'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '1.238'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +5.09%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'NEARProtocol'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.010'
$ws.Cells.Item(49, 4).ClearFormats()
$ws.Cells.Item(49, 5).Value = '  -0.20%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'Quant'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '121.36'
$ws.Cells.Item(50, 4).ClearFormats()
$ws.Cells.Item(50, 5).Value = '  -1.83%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '79.57'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  +2.46%  '
